$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while forcing Text format so that
# numeric-looking strings (leading zeros, trailing spaces, date-looking
# strings) are kept as literal text instead of being auto-converted by
# Excel into numbers / dates. Style is reset back to "Normal" afterwards
# so no stray number-format style is left behind on the cell.
function Set-TextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# --- Update existing "Pendiente ADM" OT values with real OT numbers ---
Set-TextCell "E13" "01604476 "
Set-TextCell "E23" "01549429"
Set-TextCell "E24" "01549537 "
Set-TextCell "E25" "01549613 "
Set-TextCell "E26" "01565483 "
Set-TextCell "E27" "01565494 "

# --- Append new row 28 ---
Set-TextCell "A28" "7878 "
Set-TextCell "B28" "11/26/2025"
$ws.Range("C28").Value = "LA PAMPA 3635"
$ws.Range("D28").Value = 13
Set-TextCell "E28" "01593444 "
$ws.Range("F28").Value = "Optical Power"
$ws.Range("G28").Value = "Pendiente"
$ws.Range("H28").Value = "tendido bajo"
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = '{"direccionesNormalizadas": [{"altura": 3635, "cod_calle": 12168, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.466521", "y": "-34.571932"}, "direccion": "LA PAMPA 3635, CABA", "nombre_calle": "LA PAMPA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K28").Value = -58.466521
$ws.Range("L28").Value = -34.571932
$ws.Range("M28").Value = "Colegiales"
$ws.Range("N28").Value = "Capital Norte"

# --- Append new row 29 ---
$ws.Range("A29").Value = "S00964409"
Set-TextCell "B29" "11/28/2025"
$ws.Range("C29").Value = "CERETTI 3556"
$ws.Range("D29").Value = 12
$ws.Range("E29").Value = "Pendiente ADM"
$ws.Range("F29").Value = "Optical Power"
$ws.Range("G29").Value = "Pendiente"
$ws.Range("H29").Value = "cable cortado"
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = '{"direccionesNormalizadas": [{"altura": 3556, "cod_calle": 3115, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.502145", "y": "-34.566981"}, "direccion": "CERETTI 3556, CABA", "nombre_calle": "CERETTI", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K29").Value = -58.502145
$ws.Range("L29").Value = -34.566981
$ws.Range("M29").Value = "Paternal"
$ws.Range("N29").Value = "Capital Norte"
